$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.998.46'
$ws.Range('E2').Value = '  -1.21%  '

$ws.Range('D3').Value = '3.408.61'
$ws.Range('E3').Value = '  -1.92%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = "'407.36"
$ws.Range('E5').Value = '  -1.18%  '

$ws.Range('D6').Value = "'134.26"
$ws.Range('E6').Value = '  +4.27%  '

$ws.Range('E7').Value = '  -1.20%  '

$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.15%  '

$ws.Range('D9').Value = "'0.687"
$ws.Range('E9').Value = '  -2.05%  '

$ws.Range('E10').Value = '  -6.15%  '

$ws.Range('D11').Value = "'42.72"
$ws.Range('E11').Value = '  -1.76%  '

$ws.Range('D13').Value = "'8.42"
$ws.Range('E13').Value = '  -4.03%  '

$ws.Range('D14').Value = "'19.87"
$ws.Range('E14').Value = '  -2.02%  '

$ws.Range('D15').Value = '3.424.55'
$ws.Range('E15').Value = '  -1.96%  '

$ws.Range('D16').Value = '61.986.12'
$ws.Range('E16').Value = '  -1.04%  '

$ws.Range('E17').Value = '  -3.19%  '

$ws.Range('D18').Value = "'11.02"
$ws.Range('E18').Value = '  -1.20%  '

$ws.Range('E19').Value = '  -5.79%  '

$ws.Range('E20').Value = '  -4.98%  '

$ws.Range('D21').Value = "'84.18"
$ws.Range('E21').Value = '  +2.20%  '

$ws.Range('D22').Value = "'314.65"
$ws.Range('E22').Value = '  +0.33%  '

$ws.Range('D23').Value = "'12.86"
$ws.Range('E23').Value = '  -2.85%  '

$ws.Range('D24').Value = "'3.17"
$ws.Range('E24').Value = '  -0.23%  '

$ws.Range('D25').Value = "'4.76"
$ws.Range('E25').Value = '  +9.13%  '

$ws.Range('D26').Value = "'29.58"
$ws.Range('E26').Value = '  -2.83%  '

$ws.Range('D27').Value = "'8.18"
$ws.Range('E27').Value = '  -0.04%  '

$ws.Range('D28').Value = "'2.81"
$ws.Range('E28').Value = '  +4.34%  '

$ws.Range('D29').Value = "'7.60"
$ws.Range('E29').Value = '  -2.95%  '

$ws.Range('E30').Value = '  -3.76%  '

$ws.Range('D31').Value = "'0.116"
$ws.Range('E31').Value = '  -3.84%  '

$ws.Range('D32').Value = "'42.92"
$ws.Range('E32').Value = '  -4.30%  '

$ws.Range('D33').Value = "'1.00"
$ws.Range('E33').Value = '  -0.12%  '

$ws.Range('E34').Value = '  -6.41%  '

$ws.Range('E35').Value = '  -2.73%  '

$ws.Range('D36').Value = "'51.73"
$ws.Range('E36').Value = '  -1.55%  '

$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  +0.24%  '

$ws.Range('E38').Value = '  -4.66%  '

$ws.Range('E40').Value = '  -0.66%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = "'0.125"
$ws.Range('E41').Value = '  -0.66%  '

$ws.Range('D42').Value = "'137.26"
$ws.Range('E42').Value = '  -0.37%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = "'0.299"
$ws.Range('E43').Value = '  +3.20%  '

$ws.Range('E44').Value = '  +0.81%  '

$ws.Range('D45').Value = "'16.72"
$ws.Range('E45').Value = '  -6.65%  '

$ws.Range('E46').Value = '  -2.67%  '

$ws.Range('D47').Value = "'21.36"
$ws.Range('E47').Value = '  -5.54%  '

$ws.Range('D48').Value = '2.122.20'
$ws.Range('E48').Value = '  -4.48%  '

$ws.Range('D49').Value = "'2.33"
$ws.Range('E49').Value = '  -2.50%  '

$ws.Range('D50').Value = "'1.94"
$ws.Range('E50').Value = '  +2.25%  '

$ws.Range('E51').Value = '  +16.78%  '
